# tutorial-06-migratability: hide the "Object Serialization Using PUP"
# subsection title slide (slide 2) and give it a slow slide transition.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Hide the slide from the slide show (-> <p:sld ... show="0">)
$s.SlideShowTransition.Hidden = -1

# Give the slide a "slow" transition lasting 2 seconds.
# Duration must be set before Speed so both survive in the serialized
# <p:transition/> element.
$s.SlideShowTransition.Duration = 2
$s.SlideShowTransition.Speed = 1
